$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 84: drop the trailing empty H84/I84 placeholder cells ---
# (they move down to the new last row, H88/I88, below)
$ws.Range("H84").Value = ""
$ws.Range("I84").Value = ""

# --- New rows 85-88 -------------------------------------------------
# Column A holds dates written as literal text (e.g. "2024-05-20"), not
# real date values, matching the source data. Force text entry via a
# temporary "@" number format, then reset the style back to Normal so no
# lingering custom format is left applied to the cell.
$dateCells = "A85:A88"
$ws.Range($dateCells).NumberFormat = "@"

# Row 85
$ws.Range("A85").Value = "2024-05-20"
$ws.Range("B85").Value = "13:43:58"
$ws.Range("C85").Value = "-"
$ws.Range("D85").Value = "-"
$ws.Range("E85").Value = "Power atascado en prensa, cuesta sacar"
$ws.Range("F85").Value = "-"
$ws.Range("G85").Value = "-"

# Row 86
$ws.Range("A86").Value = "2024-05-20"
$ws.Range("B86").Value = "13:44:52"
$ws.Range("C86").Value = "-"
$ws.Range("D86").Value = "-"
$ws.Range("E86").Value = "Tornillo atascado"
$ws.Range("F86").Value = "-"
$ws.Range("G86").Value = "-"

# Row 87
$ws.Range("A87").Value = "2024-05-20"
$ws.Range("B87").Value = "13:58:53"
$ws.Range("C87").Value = "Fallo en paletizador"
$ws.Range("D87").Value = "-"
$ws.Range("E87").Value = "-"
$ws.Range("F87").Value = "-"
$ws.Range("G87").Value = "-"

# Row 88
$ws.Range("A88").Value = "2024-05-20"
$ws.Range("B88").Value = "14:24:33"
$ws.Range("C88").Value = "Ascensor no sube"
$ws.Range("D88").Value = "-"
$ws.Range("E88").Value = "-"
$ws.Range("F88").Value = "-"
$ws.Range("G88").Value = "-"

# Clear the temporary format so A85:A88 end up back on the default style.
$ws.Range($dateCells).Style = "Normal"

# H88/I88: empty placeholder cells (same pattern the old H84/I84 had).
# A bare "" assignment removes the cell outright, so force a present,
# empty text cell via a quote-prefixed entry, then drop the format again.
$ws.Range("H88").NumberFormat = "@"
$ws.Range("I88").NumberFormat = "@"
$ws.Range("H88").Value = "'"
$ws.Range("I88").Value = "'"
$ws.Range("H88:I88").Style = "Normal"
